# Assignment 2 ML.docx - "Hyperparameter Tuning" commit
#
# The Q1 prompt paragraph is made bold in its entirety, and a new
# explanatory paragraph ("The following code is a python program ...")
# is inserted right after it. The _GoBack bookmark, which used to sit
# between "Regularized Logistic Regression" and ". Comment on the code
# ..." moves to the end of the newly added paragraph.

$d = $word.ActiveDocument

# Locate the "Q1. Utilize the given Jupyter notebook ..." paragraph by
# searching for its distinctive text rather than a hard-coded index.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Q1. Utilize the given Jupyter notebook*") {
        $target = $para
    }
}

# Split the paragraph right after its current end (i.e. after "...
# explaining utilized Machine Learning concepts where necessary "),
# creating a brand new, empty paragraph after it that inherits the same
# paragraph formatting (tabs/spacing/indent).
$splitPoint = $target.Range
$splitPoint.Collapse(0)
$splitPoint.InsertParagraphAfter()

# The newly created paragraph is now immediately after $target.
$newPara = $target.Next()
$newParaRange = $newPara.Range

# Insert the explanatory sentence, plus a throwaway placeholder
# character. The placeholder gives us a safe (non paragraph-final)
# anchor position to drop the relocated _GoBack bookmark at, which we
# remove immediately afterwards, leaving the bookmark collapsed right
# at the end of the inserted text (matching the target layout).
$newText = "The following code is a python program that demonstrates " + `
    "regularized logistic regression. Logistic Regression is a type " + `
    "of statistical model used to classify data into binary outcomes. " + `
    "It is a supervised learning algorithm "
$newParaRange.InsertAfter($newText + "X")

$newPara2 = $target.Next()
$bmPos = $newPara2.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$newPara3 = $target.Next()
$placeholderPos = $newPara3.Range.End - 2
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()

# Finally, bold the whole original Q1 paragraph (paragraph mark plus
# every run in it, including the sentence that used to trail the
# bookmark and now ends the paragraph).
$target.Range.Font.Bold = 1
